$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style from an existing header cell (A1) onto the new header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 114  # AD
    $ws.Cells.Item($r, 31).Value = 48   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
